# Auto-generated edit script applying the Omega_Profits value updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1740.8572
$ws.Range("J19").Value = 2037.5714
$ws.Range("L19").Value = 2037.5714
$ws.Range("N19").Value = -2387.5714
$ws.Range("H118").Value = 1236.8667
$ws.Range("I118").Value = 1236.8667
$ws.Range("K118").Value = 3710.6001
$ws.Range("M118").Value = -2053.6001
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H132").Value = 2541.309
$ws.Range("I132").Value = 2226.3845
$ws.Range("K132").Value = 6679.1535
$ws.Range("M132").Value = -4149.1535
$ws.Range("H138").Value = 5390.357
$ws.Range("J138").Value = 6988.579
$ws.Range("L138").Value = 20965.737
$ws.Range("N138").Value = -31245.737

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5200.645
$ws.Range("I61").Value = 3680.0417
$ws.Range("K61").Value = 3680.0417
$ws.Range("M61").Value = -3468.0417
$ws.Range("H101").Value = 109995
$ws.Range("J101").Value = 109995
$ws.Range("L101").Value = 109995
$ws.Range("N101").Value = -116485
$ws.Range("H102").Value = 2075.1667
$ws.Range("I102").Value = 2228.875
$ws.Range("K102").Value = 2228.875
$ws.Range("M102").Value = -606.875
$ws.Range("H107").Value = 49945
$ws.Range("J107").Value = 49945
$ws.Range("L107").Value = 49945
$ws.Range("N107").Value = -57625
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H112").Value = 27497.666
$ws.Range("J112").Value = 27497.666
$ws.Range("L112").Value = 27497.666
$ws.Range("N112").Value = -30451.666
$ws.Range("H132").Value = 4168.095
$ws.Range("I132").Value = 3922.6843
$ws.Range("K132").Value = 11768.0529
$ws.Range("M132").Value = -9238.052899999999
$ws.Range("H136").Value = 5200.645
$ws.Range("I136").Value = 3680.0417
$ws.Range("K136").Value = 11040.1251
$ws.Range("M136").Value = -8490.125100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1487.3636
$ws.Range("I107").Value = 1487.3636
$ws.Range("K107").Value = 1487.3636
$ws.Range("M107").Value = 432.6364000000001
$ws.Range("H134").Value = 3532.125
$ws.Range("I134").Value = 3532.125
$ws.Range("K134").Value = 10596.375
$ws.Range("M134").Value = -8061.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6685.6035
$ws.Range("I31").Value = 12082.857
$ws.Range("J31").Value = 4968.2954
$ws.Range("K31").Value = 12082.857
$ws.Range("L31").Value = 4968.2954
$ws.Range("M31").Value = -11787.857
$ws.Range("N31").Value = -5558.2954
$ws.Range("H34").Value = 6685.6035
$ws.Range("I34").Value = 12082.857
$ws.Range("J34").Value = 4968.2954
$ws.Range("K34").Value = 12082.857
$ws.Range("L34").Value = 4968.2954
$ws.Range("M34").Value = -11880.857
$ws.Range("N34").Value = -5372.2954
$ws.Range("H58").Value = 3254.5
$ws.Range("I58").Value = 3749.8
$ws.Range("K58").Value = 3749.8
$ws.Range("M58").Value = -3546.8
$ws.Range("H74").Value = 50314
$ws.Range("J74").Value = 50314
$ws.Range("L74").Value = 50314
$ws.Range("N74").Value = -52062
$ws.Range("H77").Value = 50314
$ws.Range("J77").Value = 50314
$ws.Range("L77").Value = 150942
$ws.Range("N77").Value = -159678
$ws.Range("H88").Value = 34768.285
$ws.Range("J88").Value = 34768.285
$ws.Range("L88").Value = 34768.285
$ws.Range("N88").Value = -35580.285
$ws.Range("H91").Value = 34768.285
$ws.Range("J91").Value = 34768.285
$ws.Range("L91").Value = 34768.285
$ws.Range("N91").Value = -37576.285
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H96").Value = 21612
$ws.Range("J96").Value = 21612
$ws.Range("L96").Value = 21612
$ws.Range("N96").Value = -27104
$ws.Range("H107").Value = 22727970
$ws.Range("J107").Value = 720.8570999999999
$ws.Range("L107").Value = 720.8570999999999
$ws.Range("N107").Value = -4560.8571
$ws.Range("H136").Value = 3254.5
$ws.Range("I136").Value = 3749.8
$ws.Range("K136").Value = 11249.4
$ws.Range("M136").Value = -8699.400000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 850
$ws.Range("I22").Value = 850
$ws.Range("K22").Value = 2550
$ws.Range("M22").Value = -2381
$ws.Range("H27").Value = 850
$ws.Range("I27").Value = 850
$ws.Range("K27").Value = 2550
$ws.Range("M27").Value = -2448
$ws.Range("H46").Value = 34638.266
$ws.Range("J46").Value = 78671.38
$ws.Range("L46").Value = 236014.14
$ws.Range("N46").Value = -236196.14
$ws.Range("H121").Value = 100645.5
$ws.Range("J121").Value = 167471.67
$ws.Range("L121").Value = 502415.01
$ws.Range("N121").Value = -505035.01

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 870.44446
$ws.Range("I107").Value = 663
$ws.Range("J107").Value = 1596.5
$ws.Range("K107").Value = 663
$ws.Range("L107").Value = 1596.5
$ws.Range("M107").Value = 1257
$ws.Range("N107").Value = -5436.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 4200
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 4200
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 4200
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -4546
$ws.Range("H104").Value = 15797
$ws.Range("J104").Value = 15797
$ws.Range("L104").Value = 15797
$ws.Range("N104").Value = -22785
$ws.Range("H110").Value = 19975
$ws.Range("J110").Value = 19975
$ws.Range("L110").Value = 19975
$ws.Range("N110").Value = -28155
$ws.Range("H132").Value = 2299.4167
$ws.Range("J132").Value = 455
$ws.Range("L132").Value = 1365
$ws.Range("N132").Value = -6425

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 23127.5
$ws.Range("I29").Value = 31505
$ws.Range("J29").Value = 14750
$ws.Range("K29").Value = 31505
$ws.Range("L29").Value = 14750
$ws.Range("M29").Value = -31215
$ws.Range("N29").Value = -15330
$ws.Range("H104").Value = 12987.5
$ws.Range("J104").Value = 12987.5
$ws.Range("L104").Value = 12987.5
$ws.Range("N104").Value = -19975.5
$ws.Range("H107").Value = 298.4
$ws.Range("I107").Value = 298.5
$ws.Range("J107").Value = 298.33334
$ws.Range("K107").Value = 895.5
$ws.Range("L107").Value = 895.0000200000001
$ws.Range("M107").Value = 1024.5
$ws.Range("N107").Value = -4735.00002
$ws.Range("H140").Value = 59806
$ws.Range("J140").Value = 59806
$ws.Range("L140").Value = 59806
$ws.Range("N140").Value = -70166
